$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1860
$ws1.Range("F5").Value = 633
$ws1.Range("F6").Value = 228

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1860
$ws4.Range("F6").Value = 633
$ws4.Range("F7").Value = 228
